$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly report date (2023-07-27, serial 45134) was added for "Vega Modelo
# de Temuco" / Maracuyá, with two quality-grade rows ("Primera" and "Segunda").
# Insert 2 new rows at the top of this product's data block (row 69), pushing
# the existing 40 data rows (old rows 69-108) down to rows 71-110.
$ws.Range("A69:A70").EntireRow.Insert()

# New row 69: Primera
$ws.Range("A69").Value = 10
$ws.Range("B69").Value = "Vega Modelo de Temuco"
$ws.Range("C69").Value = "La Araucanía"
$ws.Range("D69").Value = 45134
$ws.Range("E69").Value = 9
$ws.Range("F69").Value = "Fruta"
$ws.Range("G69").Value = 100108
$ws.Range("H69").Value = "Tropicales y subtropicales"
$ws.Range("I69").Value = 100108003
$ws.Range("J69").Value = "Maracuyá"
$ws.Range("K69").Value = "Sin especificar"
$ws.Range("L69").Value = "Primera"
$ws.Range("M69").Value = 100
$ws.Range("N69").Value = 45000
$ws.Range("O69").Value = 45000
$ws.Range("P69").Value = 45000
$ws.Range("Q69").Value = "$/caja 18 kilos"
$ws.Range("R69").Value = "Región de Arica y Parinacota"
$ws.Range("S69").Value = 2500
$ws.Range("T69").Value = 18

# New row 70: Segunda
$ws.Range("A70").Value = 10
$ws.Range("B70").Value = "Vega Modelo de Temuco"
$ws.Range("C70").Value = "La Araucanía"
$ws.Range("D70").Value = 45134
$ws.Range("E70").Value = 9
$ws.Range("F70").Value = "Fruta"
$ws.Range("G70").Value = 100108
$ws.Range("H70").Value = "Tropicales y subtropicales"
$ws.Range("I70").Value = 100108003
$ws.Range("J70").Value = "Maracuyá"
$ws.Range("K70").Value = "Sin especificar"
$ws.Range("L70").Value = "Segunda"
$ws.Range("M70").Value = 50
$ws.Range("N70").Value = 25000
$ws.Range("O70").Value = 25000
$ws.Range("P70").Value = 25000
$ws.Range("Q70").Value = "$/caja 18 kilos"
$ws.Range("R70").Value = "Región de Arica y Parinacota"
$ws.Range("S70").Value = 1389
$ws.Range("T70").Value = 18
